# edit.ps1 -- applies the "added repcount task. Also documented it" commit
# to the LogoMor Documentation.docx via the Word COM-interop object model.
#
# Strategy: locate each target paragraph with Find (on plain text, which is
# robust to run-splitting), then either:
#   - for simple in-place text tweaks: assign Range.Text / add a bookmark
#   - for paragraph restructuring (new runs, proofErr markers, new
#     paragraphs, style changes): replace the whole paragraph(s) via
#     Range.InsertXML with a literal WordProcessingML fragment, which lets
#     us control run splits / rPr / proofErr exactly like the target diff.

$d = $word.ActiveDocument

function New-Pkg([string]$innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "hides the avatar" -> "Hides the avatar", with a _GoBack bookmark
#    split after the first letter (mirrors Word's last-edit-position
#    bookmark landing where the author capitalized the sentence).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("hides the avatar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + 1
$firstChar = $d.Range($rng.Start, $splitPos)
$firstChar.Text = "H"
$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

# ---------------------------------------------------------------------
# 2) Program Flow Control section becomes Conditionals / If, Ifelse ;
#    a new "Loops" H2 is introduced and the old "Loops" H3 becomes
#    "Repeat".
# ---------------------------------------------------------------------

# 2a. Heading3 "Conditionals " -> "If, " + Ifelse (spell-checked run).
#     Do this BEFORE renaming the Heading2 below it, because once that
#     rename also produces a "Conditionals " paragraph, a plain text
#     search for "Conditionals " would become ambiguous.
$rng = $d.Content
$rng.Find.Execute("Conditionals ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">If, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Ifelse</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
)
$rng.InsertXML($xml)

# 2b. Heading2 "Program Flow Control" -> "Conditionals "
$rng = $d.Content
$rng.Find.Execute("Program Flow Control", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Conditionals "

# 2c. Insert a new Heading2 "Loops" right before the existing Heading3
#     "Loops" paragraph, and rename that Heading3 to "Repeat".
$rng = $d.Content
$rng.Find.Execute("Loops", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Loops</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Repeat</w:t></w:r></w:p>'
)
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) After the "repeat" example, document the new repcount command
#    (Note: heading + paragraph), then add the "While" H3 heading.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("while :n > 0 [ fd 10 rt 90 make", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# rewind to before this search started: find the repeat-example paragraph
$rng = $d.Content
$rng.Find.Execute("repeat 4 [ fd 10 rt 90 ]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.Expand(4) | Out-Null   # the blank NoSpacing paragraph right after the example
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:r><w:t>Note:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Using the command </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>repcount</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>we can get the number of the current execution of the innermost repeat block being executed, starting from number 1. If no repeat block is executed, 0 is returned</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>While</w:t></w:r></w:p>'
)
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) After the "while" example, add the "Until" H3 heading.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("while :n > 0 [ fd 10 rt 90 make", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.Expand(4) | Out-Null   # the blank NoSpacing paragraph right after the while example
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Until</w:t></w:r></w:p>'
)
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 5) Mouse commands intro paragraph: drop the stray _GoBack bookmark and
#    merge the two runs back into a single run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The following commands provide info about the user", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$apostrophe = [char]8217
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t>The following commands provide info about the user' + $apostrophe + 's mouse position and button.</w:t></w:r>' +
    '</w:p>'
)
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) Comments example: "Fd 20 ;move a bit forward" -> "fd 20 ;move a bit
#    forward" fully in Courier New 10pt, with tightened proofErr spans.
#    We include the preceding "Example" heading paragraph in the replaced
#    range so the range doesn't start exactly on the pre-existing
#    zero-width <w:proofErr/> (which would otherwise survive untouched).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("it appeared", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null
$rng.Collapse(0) | Out-Null
$exampleStart = $rng.Start
$rng.Expand(4) | Out-Null      # "Example:" Heading5 paragraph
$rng.Collapse(0) | Out-Null
$rng.Expand(4) | Out-Null      # "Fd 20 ;move a bit forward" paragraph
$fdEnd = $rng.End

$full = $d.Range($exampleStart, $fdEnd)
$courierRPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'
$xml = New-Pkg(
    '<w:p><w:pPr><w:pStyle w:val="Heading5"/></w:pPr><w:r><w:t>Example</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $courierRPr + '<w:t>f</w:t></w:r>' +
    '<w:r>' + $courierRPr + '<w:t>d</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $courierRPr + '<w:t xml:space="preserve"> 20 ;move a bit forward</w:t></w:r>' +
    '</w:p>'
)
$full.InsertXML($xml)
